# edit suite à meeting
# Remove the whole "*lien avec la plateforme de devis ?" question paragraph
# from the QUESTIONS/SUGGESTIONS list (it's no longer relevant after the
# meeting).

$d = $word.ActiveDocument

foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*lien avec la plateforme de devis*") {
        $p.Range.Delete()
        break
    }
}
